# Coaching Spreadsheet - "Update - bug fixes"
#
# 1. Rename "Sheet1" -> "Enrollees"
# 2. Widen column B (16.28 -> 21.02 OOXML "characters") and nudge column E
#    (36.2 -> 36.19 OOXML "characters").
#
# Note: this COM runtime re-derives the stored OOXML column width from the
# `ColumnWidth` value using a fixed 6-pixel "character" grid with 5px of
# padding (stored = round(ColumnWidth*6)/6 + 5/6), so only multiples of
# 1/6 are representable. The values below are chosen so the saved width
# lands as close as possible to the target (21.02 -> 21, 36.19 -> 36.1667).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Enrollees"

$ws.Columns.Item(2).ColumnWidth = 20.166666666666668
$ws.Columns.Item(5).ColumnWidth = 35.333333333333336
